# Generate and save output file after processing
#
# Inserts three new "general_college_subjects" columns (history, electives,
# cs) right before the existing "arts" column, shifting every column from
# the old R (arts) through AE (act scores.75th) three places to the right
# (new range becomes U:AH). The insert naturally carries over the values
# and types of the shifted cells, which is exactly what the target diff
# shows. We then populate the three new header cells and the new data-row
# cells, and normalize a handful of data-row text values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 columns at R:T -- shifts old R:AE (arts .. act scores.75th)
# to U:AH, expanding the used range from A1:AE2 to A1:AH2.
$ws.Range("R1:T1").EntireColumn.Insert()

# New column headers (row 1)
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New column values for the Stanford data row (row 2)
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Normalize existing data-row text values
$ws.Range("D2").Value = "considered"
$ws.Range("E2").Value = "considered"
$ws.Range("F2").Value = "not considered"
$ws.Range("G2").Value = "considered"
$ws.Range("H2").Value = "very important"
$ws.Range("I2").Value = "considered"
$ws.Range("J2").Value = "considered"
